# Add column C ('t+3') to the predicted variables matrix, following the
# same layout as columns A and B: a header value in row 1, then one
# predicted value per row down to row 67.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label for the new column (A1=0, B1=1, C1=2 -> t+1/t+2/t+3 offsets)
$header = 2

# Predicted values for C2:C67, in row order
$values = @(
    -0.3093107257104907,
    0.2076865250416704,
    0.1425296526411839,
    -0.07633509935382571,
    0.1418282182027905,
    -0.4528029662052117,
    -0.260071556078068,
    -0.4090605937397951,
    0.358506294481997,
    -0.2401966792289666,
    -0.08648358856128772,
    -0.0135888336771998,
    0.05244956046941167,
    -0.01300978491836207,
    0.3511649265855035,
    0.6020326650123404,
    0.03866759896343645,
    0.3942978762976164,
    0.2970493814345418,
    0.5205383953746511,
    0.3797409485700736,
    -0.06102057954108653,
    4.600142545217028,
    0.4798567206283482,
    0.3525363144340392,
    0.310901809741786,
    1.066890520460006,
    5.468790468361631,
    0.9649696633972918,
    -0.2243439863291791,
    0.7474304657757006,
    0.8993922506230235,
    -0.6672725585690531,
    0.8131609638990331,
    0.767376525621137,
    0.7413393640081756,
    0.7166327809518839,
    0.5805097688759314,
    0.7508300768615521,
    0.5545927399189651,
    0.6890526931284596,
    0.708607646620943,
    0.6736650154960108,
    0.6573158281272644,
    -1.265986458938779,
    -0.9814853793243534,
    -0.8755392624003318,
    -0.6439722445150895,
    -0.05022248240209715,
    -0.8664839627003759,
    -0.8664839627003759,
    -1.099638969396312,
    -0.1879309031566863,
    -1.001920247713072,
    -0.8865783913471785,
    -0.9350879674765653,
    -1.129367752215911,
    -0.8481378087766493,
    -0.4975904977737807,
    0.370916657261192,
    -1.212315704557025,
    -0.7377023705146881,
    -0.9026255638979531,
    -0.1111563709664542,
    -0.7926290785388965,
    -0.7645593498930221
)

# Copy the formatting of the existing header cell (B1) onto the new
# header cell (C1) so it reuses the same bold/centered/bordered style.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").Value = $header

# Loop over the data rows and fill column C
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
